$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.036.27"
$ws.Range("E2").Value = "  +0.79%  "

# Row 3
$ws.Range("D3").Value = "1.904.50"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "0.8298"
$ws.Range("E5").Value = "  +9.00%  "

# Row 6
$ws.Range("D6").Value = "241.91"
$ws.Range("E6").Value = "  +0.78%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "0.3231"
$ws.Range("E8").Value = "  +6.45%  "

# Row 9
$ws.Range("D9").Value = "26.76"
$ws.Range("E9").Value = "  +5.40%  "

# Row 10
$ws.Range("D10").Value = "0.07029"
$ws.Range("E10").Value = "  +3.16%  "

# Row 11
$ws.Range("D11").Value = "0.08029"
$ws.Range("E11").Value = "  +0.72%  "

# Row 12
$ws.Range("D12").Value = "0.7507"
$ws.Range("E12").Value = "  +2.13%  "

# Row 13
$ws.Range("D13").Value = "1.905.39"
$ws.Range("E13").Value = "  +0.57%  "

# Row 14
$ws.Range("E14").Value = "  +1.31%  "

# Row 15
$ws.Range("D15").Value = "92.76"
$ws.Range("E15").Value = "  +2.16%  "

# Row 16
$ws.Range("D16").Value = "30.048.19"
$ws.Range("E16").Value = "  +0.87%  "

# Row 17
$ws.Range("D17").Value = "14.15"
$ws.Range("E17").Value = "  +2.65%  "

# Row 18
$ws.Range("D18").Value = "5.939"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19
$ws.Range("D19").Value = "244.38"
$ws.Range("E19").Value = "  +1.35%  "

# Row 20
$ws.Range("D20").Value = "0.000007778"
$ws.Range("E20").Value = "  +1.15%  "

# Row 21
$ws.Range("D21").Value = "2.162.21"
$ws.Range("E21").Value = "  +1.06%  "

# Row 22
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.08%  "

# Row 24
$ws.Range("D24").Value = "6.994"
$ws.Range("E24").Value = "  +1.51%  "

# Row 25
$ws.Range("D25").Value = "0.1580"
$ws.Range("E25").Value = "  +23.10%  "

# Row 26
$ws.Range("D26").Value = "168.49"
$ws.Range("E26").Value = "  +1.19%  "

# Row 27
$ws.Range("D27").Value = "9.235"
$ws.Range("E27").Value = "  +0.29%  "

# Row 28
$ws.Range("E28").Value = "  +2.00%  "

# Row 29
$ws.Range("E29").Value = "  +3.63%  "

# Row 30
$ws.Range("D30").Value = "1.374"
$ws.Range("E30").Value = "  -1.98%  "

# Row 31
$ws.Range("E31").Value = "  +0.24%  "

# Row 32
$ws.Range("D32").Value = "4.286"
$ws.Range("E32").Value = "  +0.80%  "

# Row 33
$ws.Range("D33").Value = "0.05714"
$ws.Range("E33").Value = "  +10.01%  "

# Row 34
$ws.Range("D34").Value = "4.095"
$ws.Range("E34").Value = "  +0.87%  "

# Row 35
$ws.Range("D35").Value = "1.289"
$ws.Range("E35").Value = "  +3.63%  "

# Row 36
$ws.Range("D36").Value = "0.7346"

# Row 37
$ws.Range("D37").Value = "2.727"
$ws.Range("E37").Value = "  +0.50%  "

# Row 39
$ws.Range("E39").Value = "  +0.84%  "

# Row 40
$ws.Range("D40").Value = "0.4423"
$ws.Range("E40").Value = "  +0.77%  "

# Row 41
$ws.Range("D41").Value = "72.31"
$ws.Range("E41").Value = "  +0.70%  "

# Row 42
$ws.Range("D42").Value = "5.964"
$ws.Range("E42").Value = "  -2.82%  "

# Row 43
$ws.Range("D43").Value = "0.8433"
$ws.Range("E43").Value = "  +1.88%  "

# Row 44
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("D45").Value = "1.897"
$ws.Range("E45").Value = "  +1.08%  "

# Row 46
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "101.39"
$ws.Range("E46").Value = "  +1.89%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.623"
$ws.Range("E47").Value = "  +0.34%  "

# Row 48
$ws.Range("D48").Value = "9.725"
$ws.Range("E48").Value = "  +0.60%  "

# Row 49
$ws.Range("D49").Value = "994.34"
$ws.Range("E49").Value = "  +9.58%  "

# Row 50
$ws.Range("D50").Value = "2.063.40"
$ws.Range("E50").Value = "  +1.17%  "

# Row 51
$ws.Range("D51").Value = "36.30"
$ws.Range("E51").Value = "  +0.76%  "
